$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells whose new value looks like a plain
# number (e.g. "0.9987" or "158.00") to stay plain text, matching the
# source workbook's inline-string cells -- otherwise Excel would
# auto-detect them as numbers and change their cell type. Cells whose
# new value already contains a second "." (e.g. "29.332.44") are never
# auto-converted by Excel, so they are left alone.
$priceCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D44", "D45", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    'D2' = '29.332.44'
    'E2' = '  -0.05%  '
    'D3' = '1.839.08'
    'E3' = '  -0.31%  '
    'D4' = '0.9987'
    'E4' = '  +0.13%  '
    'D5' = '239.13'
    'E5' = '  -0.34%  '
    'D6' = '0.6273'
    'E6' = '  +0.18%  '
    'D8' = '0.07429'
    'E8' = '  -0.82%  '
    'D9' = '0.2896'
    'E9' = '  -0.09%  '
    'D10' = '24.94'
    'E10' = '  +2.35%  '
    'D11' = '0.07723'
    'E11' = '  +0.08%  '
    'D12' = '1.846.99'
    'E12' = '  +0.13%  '
    'D13' = '4.959'
    'E13' = '  -0.73%  '
    'D14' = '0.6748'
    'E14' = '  -0.57%  '
    'D15' = '0.00001021'
    'E15' = '  -1.13%  '
    'D16' = '81.72'
    'E16' = '  -0.43%  '
    'D17' = '6.231'
    'E17' = '  +1.45%  '
    'D18' = '29.319.25'
    'E18' = '  -0.27%  '
    'D19' = '230.39'
    'E19' = '  +0.61%  '
    'D20' = '12.29'
    'E20' = '  -0.31%  '
    'D21' = '0.9998'
    'E21' = '  +0.12%  '
    'D22' = '7.346'
    'E22' = '  -1.55%  '
    'D23' = '0.9999'
    'E23' = '  +0.22%  '
    'D24' = '158.00'
    'E24' = '  -0.30%  '
    'D25' = '8.481'
    'E25' = '  +0.95%  '
    'E26' = '  -1.89%  '
    'D27' = '17.35'
    'E27' = '  -0.92%  '
    'D28' = '0.07262'
    'E28' = '  +12.28%  '
    'D29' = '1.456'
    'E29' = '  +5.48%  '
    'D30' = '1.477'
    'E30' = '  +0.36%  '
    'D31' = '4.040'
    'E31' = '  -1.23%  '
    'D32' = '4.041'
    'E32' = '  -0.54%  '
    'D33' = '1.816'
    'E33' = '  -0.45%  '
    'D34' = '1.140'
    'E34' = '  +0.00%  '
    'D35' = '0.6953'
    'E35' = '  -0.64%  '
    'D36' = '2.570'
    'E36' = '  -0.21%  '
    'D37' = '0.01836'
    'E37' = '  +0.52%  '
    'D38' = '2.810'
    'E38' = '  -0.83%  '
    'D39' = '6.847'
    'E39' = '  +3.66%  '
    'D40' = '1.234.04'
    'E40' = '  -2.03%  '
    'D41' = '0.9361'
    'E41' = '  +2.83%  '
    'D42' = '1.000'
    'E42' = '  +0.19%  '
    'D43' = '1.986.83'
    'E43' = '  -0.96%  '
    'D44' = '100.97'
    'E44' = '  -0.50%  '
    'D45' = '65.39'
    'E45' = '  -1.06%  '
    'E46' = '  +0.85%  '
    'D47' = '1.704'
    'E47' = '  -1.51%  '
    'D48' = '6.944'
    'E48' = '  -1.78%  '
    'D49' = '0.1139'
    'E49' = '  -2.88%  '
    'D50' = '8.893'
    'E50' = '  -1.23%  '
    'D51' = '0.3901'
    'E51' = '  -1.01%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Host "Updated $($updates.Count) cells"
